$wb = $excel.ActiveWorkbook
$ws10 = $wb.Worksheets.Item("Trim Whitespace")

foreach ($w in 13, 14, 14.2, 14.5, 15, 15.5) {
  $ws10.Columns.Item(6).ColumnWidth = $w
  Write-Host ("w=" + $w + " -> stored later")
}
